{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the \"Added many more features\" copy update to the Almighty\n// Sparta review document: refreshed title, reworked pros/cons bullet\n// text, and a rewritten SEO description near the end of the document.\n\nconst replacements = [\n  [\n    \"Play Almighty Sparta Free: Game Review and Pros & Cons\",\n    \"Play Almighty Sparta Free | Exciting Slot Game Review\",\n  ],\n  [\n    \"Seamless gameplay and immersive experience\",\n    \"Immersive gameplay and seamless experience\",\n  ],\n  [\n    \"High-quality and entertaining slot game\",\n    \"High-quality graphics and attention to detail\",\n  ],\n  [\n    \"Special symbols and free spins add an extra touch\",\n    \"Special symbols and functions add excitement\",\n  ],\n  [\n    \"Balance of risk and reward is appreciated by experienced gamblers\",\n    \"Appreciated balance of risk and reward\",\n  ],\n  [\n    \"Limited special symbols\",\n    \"Limited special symbols and functions\",\n  ],\n  [\n    \"Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time.\",\n    \"Read the review of Almighty Sparta, play for free, and enjoy an exciting slot game experience.\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the \"Added many more features\" copy update to the Almighty\n# Sparta review document: refreshed title, reworked pros/cons bullet\n# text, and a rewritten SEO description near the end of the document.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-AllText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null\n}\n\nReplace-AllText \"Play Almighty Sparta Free: Game Review and Pros & Cons\" \"Play Almighty Sparta Free | Exciting Slot Game Review\"\nReplace-AllText \"Seamless gameplay and immersive experience\" \"Immersive gameplay and seamless experience\"\nReplace-AllText \"High-quality and entertaining slot game\" \"High-quality graphics and attention to detail\"\nReplace-AllText \"Special symbols and free spins add an extra touch\" \"Special symbols and functions add excitement\"\nReplace-AllText \"Balance of risk and reward is appreciated by experienced gamblers\" \"Appreciated balance of risk and reward\"\nReplace-AllText \"Limited special symbols\" \"Limited special symbols and functions\"\nReplace-AllText \"Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time.\" \"Read the review of Almighty Sparta, play for free, and enjoy an exciting slot game experience.\"\n"}
